# The "Förändrad" (Changed) date in column C was updated from 2023-09-03
# (serial 45172) to 2023-09-06 (serial 45175) for every data row (rows 2-422).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C422").Value = 45175
